# Update the 12 "Salm" QC data rows (rows 2-13) on the "Data" sheet to the
# next run's values (20220928-Salm-14911Updt -> 20221012-Salm-10004Updt):
#   - Result ID (A): A1013601-A1013612  -> A1988301-A1988312 (row order preserved)
#   - Lab Sample ID (E): 20220928-Salm-14911Updt -> 20221012-Salm-10004Updt
#   - Lane (Q): re-sequenced to 1..12 in row order
#   - Cartridge ID (T): CartridgeSalm4911 -> CartridgeSalm0004

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $lane = $i + 1

    $resultId = "A19883{0:D2}" -f $lane

    $ws.Range("A$row").Value = $resultId
    $ws.Range("E$row").Value = "20221012-Salm-10004Updt"

    # Lane is stored as text (shared string), not a number, in the source
    # workbook, so force text formatting before assigning a numeric-looking
    # string, then restore the default style so no visible formatting change
    # is introduced.
    $ws.Range("Q$row").NumberFormat = "@"
    $ws.Range("Q$row").Value = [string]$lane
    $ws.Range("Q$row").Style = "Normal"

    $ws.Range("T$row").Value = "CartridgeSalm0004"
}
